# Applies the "may 9th" data update to the struggle sensor-log sheet:
# 5 new rows of samples are inserted at the top of the series (pushing the
# existing rows down by 5 timestamps-worth) and 5 new rows are appended at
# the end, extending the used range from A1:H21 to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number (2..31) -> full row values (A:timestamp, B:label, C..H:sensor axes)
$rows = @{
    2 = @(0, "struggle", 3.055751519543784, -7.220968450818743, -2.894419597727912, -0.09269879758358, 0.0145080499351024, -0.007177666760981)
    3 = @(100, "struggle", 2.955562557492938, -7.134888444628034, -2.944399050303869, -0.090408056974411, -0.0394008085131645, -0.0087048299610614)
    4 = @(200, "struggle", 2.833344757556915, -7.338198423385621, -2.509933024644852, -0.0737619772553443, -0.102472648024559, 0.0461203269660472)
    5 = @(300, "struggle", 3.108331612178258, -7.027578847748893, -2.627250722476415, -0.1244637966156005, -0.4952589869499206, 0.25641068816185)
    6 = @(400, "struggle", 3.771172387259348, -7.547038768018995, -3.545775217669351, 0.955545961856842, 0.418595403432846, 0.5012149214744568)
    7 = @(500, "struggle", 2.675201347896034, -6.581202302660263, 2.435893160956232, 0.2005165219306945, 3.583182811737061, 0.4506658315658569)
    8 = @(600, "struggle", 0.8038081441606844, -7.229396700859067, 1.553160754697669, 0.0250454749912023, 1.217149019241333, 0.2982549667358398)
    9 = @(700, "struggle", 1.023633156503948, -4.972137170178557, -3.916851524795791, 0.4430300295352936, 1.049771904945374, -0.152105450630188)
    10 = @(800, "struggle", 0.3411836709295031, -4.24002621429307, -5.904611808913096, -0.1605048477649688, 0.6899722814559937, -0.2005165219306945)
    11 = @(900, "struggle", -0.6227953433990485, -4.929600954055787, -6.7005569934845, 0.2364048510789871, -1.487151384353638, -0.1020144969224929)
    12 = @(1000, "struggle", -0.2611346755709001, -5.176580480166844, -11.08990749291009, -1.693776607513428, -5.596747398376465, 1.112385630607605)
    13 = @(1100, "struggle", 7.22266825607841, -7.398067985262175, -7.469464063644459, -1.431409955024719, -1.327410221099854, 1.870163917541504)
    14 = @(1200, "struggle", 17.43506016050075, -11.15922947440832, 8.020842296736639, 0.0817032232880592, -2.866485118865967, 0.2145664244890213)
    15 = @(1300, "struggle", 4.810774190085277, -6.605375971112934, 0.3875993319920126, 0.2112066596746444, -2.239890098571777, -0.6490443348884583)
    16 = @(1400, "struggle", 3.291332300220217, -5.354068347385952, 1.089433806283133, -0.4173736870288849, -0.5337435007095337, -0.0603229440748691)
    17 = @(1500, "struggle", 3.138020081179482, -6.129925046648298, 3.955937453678676, 0.2202169150114059, 1.032209515571594, 0.3608686327934265)
    18 = @(1600, "struggle", 3.258425533771516, -4.712551474571227, 5.692880451679232, 2.595261096954346, 7.212944030761719, -0.4940372705459595)
    19 = @(1700, "struggle", -3.766383392470183, -8.431432792118596, 7.765124661581846, 0.2443460971117019, 2.484389066696167, 0.5590944290161133)
    20 = @(1800, "struggle", 1.177985225405012, -3.925881973334721, -4.60967251232692, 0.3139847218990326, 0.8843801617622375, -0.1221730485558509)
    21 = @(1900, "struggle", 1.052592243467076, -3.166507703917357, -6.040607401302906, -1.060767531394958, 1.06214189529419, 0.0241291765123605)
    22 = @(2000, "struggle", -1.074721106461116, -3.308761754206249, -6.224646031856537, -0.1357648074626922, -0.5329799056053162, -0.2561052441596985)
    23 = @(2100, "struggle", -1.158144678388323, -3.9557591165815, -8.491945947919575, -0.5283984541893005, -2.495231866836548, 0.1453859210014343)
    24 = @(2200, "struggle", 0.7903139420917948, -5.482765521321981, -6.279002168348867, -0.8752171993255615, -1.659415483474731, -0.0487165041267871)
    25 = @(2300, "struggle", 4.010212659835815, -5.284669637680054, -2.184034883975983, -0.0806342139840126, -0.5641340613365173, 0.087353728711605)
    26 = @(2400, "struggle", 2.293728096144541, -7.460245260170529, 0.3245020040443973, -0.1470658034086227, -0.0910189226269722, -0.2229658216238021)
    27 = @(2500, "struggle", 1.675990547452653, -6.526311159133912, -0.6771522419793252, 0.4641048610210418, 0.3608686327934265, 0.3602577745914459)
    28 = @(2600, "struggle", 3.116939672401965, -5.273013770580299, -3.711685695818474, 0.2412917762994766, 0.2144137024879455, -0.0186313893646001)
    29 = @(2700, "struggle", 2.304866858891079, -5.722443546567645, -4.230944650513785, 0.1291979998350143, 0.1751656085252761, -0.070249505341053)
    30 = @(2800, "struggle", 2.380605901990617, -5.519177751881738, -4.144632569381169, -0.1421788930892944, -0.0774271711707115, -0.0045814891345798)
    31 = @(2900, "struggle", 2.426761286599295, -5.405185426984516, -3.328171287264147, -0.067195177078247, -0.219300627708435, -0.2028072625398636)
}

foreach ($r in $rows.Keys) {
    $rowValues = $rows[$r]
    $arr = New-Object "object[,]" 1,8
    for ($c = 0; $c -lt 8; $c++) { $arr[0, $c] = $rowValues[$c] }
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 8))
    $rng.Value = $arr
}

Write-Host "Updated rows 2-31 ($($rows.Count) data rows); used range now $($ws.UsedRange.Address())"
